# Edit summary:
#  1. The stray "_GoBack" bookmark that used to sit at the end of the
#     document title ("... Azure SQL Data Warehouse") is removed - it is
#     replaced by a fresh "_GoBack" bookmark marking the location of the
#     actual text edit below (Word keeps only one "_GoBack" bookmark per
#     document, so adding the new one automatically drops the old one).
#  2. The paragraph that reads "To avoid unnecessary cost, ake sure you
#     pause ..." is missing the leading "m" of "make" (and still carries
#     the spell-check proofErr markers Word had drawn around the
#     misspelled fragment "ake"). We type the missing "m" back in,
#     clear out the now-stale proofErr markers, and leave the new
#     "_GoBack" bookmark sitting right after the "m", matching where the
#     edit actually happened.

$d = $word.ActiveDocument

# --- locate the edit point ------------------------------------------------
$find = $d.Content
[void]$find.Find.Execute("To avoid unnecessary cost, ")
$insPt = $find.End

# --- 1. type the missing "m" as its own run -------------------------------
$ins = $d.Range($insPt, $insPt)
$ins.InsertAfter("m")

# Touch formatting so the engine keeps "m" as a distinct run instead of
# silently re-merging it into the identically-formatted preceding run.
$mRun = $d.Range($insPt, $insPt + 1)
$mRun.Font.Bold = 1
$mRun.Font.Bold = 0

$afterM = $insPt + 1

# --- 2. clear the stale spell-check markers around the old "ake" ----------
# Re-typing "make sure" through Find/Replace drops the now-irrelevant
# proofErr spellStart/spellEnd pair Word had drawn around "ake", the same
# way Word's own proofing pass removes them once the word is spelled
# correctly.
$clean = $d.Range($afterM - 1, $afterM + 5)
[void]$clean.Find.Execute($clean.Text, $true, $false, $false, $false, $false, `
                           $true, 1, $false, $clean.Text, 2)

# The replace above can re-merge "m" back into the preceding run; split it
# back out again.
$mRun2 = $d.Range($insPt, $insPt + 1)
$mRun2.Font.Bold = 1
$mRun2.Font.Bold = 0

# --- 3. drop the "_GoBack" bookmark here (replaces the old one) ----------
$bmRange = $d.Range($afterM, $afterM)
$d.Bookmarks.Add("_GoBack", $bmRange)
